# Add more scraped data about a player: a new "Player Info" sheet placed
# before the existing sheets, and convert the MATCH_CARD_LINK columns on
# the batting/bowling sheets into plain MATCH_CODE values (the bare numeric
# id that used to be the query string of the scorecard URL).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet as the first tab.
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Mirror the bold / bordered / centered header style used by the other sheets.
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").Borders.LineStyle = 1
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160

$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4331"
$playerInfo.Range("B2").Value = "Sean A Abbott"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ------------------------------------------------------------------
# 2. ODI Batting: rename MATCH_CARD_LINK -> MATCH_CODE, and replace the
#    full scorecard URL in column D with just the numeric match code.
#    (Re-fetch the sheet reference: inserting a sheet above invalidates
#    the handle grabbed earlier.)
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingRows = $battingSheet.UsedRange.Rows.Count
$battingSheet.Range("D2:D" + $battingRows).NumberFormat = "@"
for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value()
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.Value = $code
    }
}

# ------------------------------------------------------------------
# 3. ODI Bowling: same MATCH_CARD_LINK -> MATCH_CODE treatment, but the
#    link lives in column B on this sheet.
# ------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingRows = $bowlingSheet.UsedRange.Rows.Count
$bowlingSheet.Range("B2:B" + $bowlingRows).NumberFormat = "@"
for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value()
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.Value = $code
    }
}
